# Apply the authored changes to the "source" worksheet:
#   1. Move the active selection from H5 to C4.
#   2. Widen column A (was ~26.86 chars) to ~46.43 chars.
#
# (The workbook.xml metadata touched by the original diff -- fileVersion
#  rupBuild, the x15ac:absPath, the xr:revisionPtr GUID/coauth versions and
#  the bookViews window geometry -- are save-environment/machine specific
#  values that Excel stamps in from the actual host session. They are not
#  exposed as settable properties on the documented Workbook/Window COM
#  object model, so they are intentionally left alone here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Selection: H5 -> C4
$ws.Range("C4").Select()

# 2) Column A width: 26.85546875 -> 46.42578125 (stored "width" units).
#    The ColumnWidth COM property is quantized by this host to 1/6-character
#    steps, so we pick the value whose rounded result lands on the closest
#    achievable stored width (46.5) to the target 46.42578125.
$ws.Columns.Item(1).ColumnWidth = 45.65
